# Update commonness index calculations (H and I columns, rows 2-46)
# in the betas_traits_noeffort table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ H = "0.12436080724000931"; I = "0.52366214990615845" }
    3 = @{ H = "0.10100854188203812"; I = "0.50216466188430786" }
    4 = @{ H = "0.083859793841838837"; I = "0.49338066577911377" }
    5 = @{ H = "0.76485836505889893"; I = "0.9034649133682251" }
    6 = @{ H = "0.76121127605438232"; I = "0.89929115772247314" }
    7 = @{ H = "0.77029860019683838"; I = "0.90083026885986328" }
    8 = @{ H = "0.22846108675003052"; I = "1.1147490739822388" }
    9 = @{ H = "0.27701675891876221"; I = "1.10173499584198" }
    10 = @{ H = "0.28289857506752014"; I = "1.1244121789932251" }
    11 = @{ H = "0.060346517711877823"; I = "0.47596254944801331" }
    12 = @{ H = "0.026914555579423904"; I = "0.45807254314422607" }
    13 = @{ H = "0.020564153790473938"; I = "0.45535069704055786" }
    14 = @{ H = "-0.45975089073181152"; I = "0.26661381125450134" }
    15 = @{ H = "-0.43129760026931763"; I = "0.23487792909145355" }
    16 = @{ H = "-0.45413625240325928"; I = "0.26270231604576111" }
    17 = @{ H = "-0.018348162993788719"; I = "0.14830964803695679" }
    18 = @{ H = "-0.014652900397777557"; I = "0.1413377970457077" }
    19 = @{ H = "-0.015011650510132313"; I = "0.13853439688682556" }
    20 = @{ H = "0.21055229008197784"; I = "0.30660146474838257" }
    21 = @{ H = "0.20894399285316467"; I = "0.30526670813560486" }
    22 = @{ H = "0.21240603923797607"; I = "0.30484408140182495" }
    23 = @{ H = "4.1820697784423828"; I = "6.0359458923339844" }
    24 = @{ H = "4.1832270622253418"; I = "6.0716557502746582" }
    25 = @{ H = "4.1914629936218262"; I = "6.0218243598937988" }
    26 = @{ H = "-0.0078555736690759659"; I = "0.042446907609701157" }
    27 = @{ H = "-0.0089985299855470657"; I = "0.040426496416330338" }
    28 = @{ H = "-0.0086145224049687386"; I = "0.041145749390125275" }
    29 = @{ H = "-0.18132717907428741"; I = "0.80455845594406128" }
    30 = @{ H = "-0.17077110707759857"; I = "0.77271348237991333" }
    31 = @{ H = "-0.16830363869667053"; I = "0.72662299871444702" }
    32 = @{ H = "0.020239254459738731"; I = "0.11567499488592148" }
    33 = @{ H = "0.016656436026096344"; I = "0.11375459283590317" }
    34 = @{ H = "0.01699206605553627"; I = "0.11203553527593613" }
    35 = @{ H = "0.20910288393497467"; I = "0.38475227355957031" }
    36 = @{ H = "0.20453935861587524"; I = "0.38774794340133667" }
    37 = @{ H = "0.20554928481578827"; I = "0.38808223605155945" }
    38 = @{ H = "4.0995664596557617"; I = "6.6781644821166992" }
    39 = @{ H = "4.076085090637207"; I = "6.7700862884521484" }
    40 = @{ H = "4.0829825401306152"; I = "6.7392764091491699" }
    41 = @{ H = "0.0014932226622477174"; I = "0.043444715440273285" }
    42 = @{ H = "2.8386413760017604e-05"; I = "0.041798625141382217" }
    43 = @{ H = "0.0008566506439819932"; I = "0.041462104767560959" }
    44 = @{ H = "0.028615837916731834"; I = "0.55835682153701782" }
    45 = @{ H = "0.013285844586789608"; I = "0.54966199398040771" }
    46 = @{ H = "0.011353068985044956"; I = "0.53180825710296631" }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("H$row").Value = [double]$vals.H
    $ws.Range("I$row").Value = [double]$vals.I
}

$wb.Save()
